# 2026-02-22 12:40 OLX monitor update.
# Appends 8 freshly-scraped listing rows (139-146) to the bottom of the
# "PODSUMOWANIE" sheet listing table (columns A-H), matching the layout and
# formatting already used by the existing listing rows (7-138):
#   A "checked at", B "profile", C "title", D "price", E "date added",
#   F "days listed", G "url", H "slug".
# Column styling already on the sheet:
#   A, C  -> left-aligned plain text
#   D, E  -> centered
#   F     -> centered; red font once the listing has been up > 60 days
#   B, G, H -> default / unstyled

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PODSUMOWANIE")

# ---- Row 139 ----
$ws.Range("A139").Value = '2026-02-22 12:40:56'
$ws.Range("B139").Value = 'poqui'
$ws.Range("C139").Value = 'Mieszkanie z KLIMATYZACJĄ 5 minut od UMCS, UP, KUL - Długosza'
$ws.Range("D139").Value = 2049
$ws.Range("E139").Value = '19.12.2025'
$ws.Range("F139").Value = 64
$ws.Range("G139").Value = 'https://www.olx.pl/d/oferta/mieszkanie-z-klimatyzacja-5-minut-od-umcs-up-kul-dlugosza-CID3-ID18KAEc.html'
$ws.Range("H139").Value = 'mieszkanie-z-klimatyzacja-5-minut-od-umcs-up-kul-dlugosza-CID3-ID18KAEc'

# ---- Row 140 ----
$ws.Range("A140").Value = '2026-02-22 12:40:56'
$ws.Range("B140").Value = 'poqui'
$ws.Range("C140").Value = 'Świeżo wykończone mieszkanie z dużym balkonem - Ponikwoda'
$ws.Range("D140").Value = 2299
$ws.Range("E140").Value = '19.01.2026'
$ws.Range("F140").Value = 34
$ws.Range("G140").Value = 'https://www.olx.pl/d/oferta/swiezo-wykonczone-mieszkanie-z-duzym-balkonem-ponikwoda-CID3-ID1951OR.html'
$ws.Range("H140").Value = 'swiezo-wykonczone-mieszkanie-z-duzym-balkonem-ponikwoda-CID3-ID1951OR'

# ---- Row 141 ----
$ws.Range("A141").Value = '2026-02-22 12:40:56'
$ws.Range("B141").Value = 'poqui'
$ws.Range("C141").Value = 'Kawalerka po remoncie z funkcjonalną antresolą - ul. Jana Sawy'
$ws.Range("D141").Value = 2499
$ws.Range("E141").Value = '28.10.2025'
$ws.Range("F141").Value = 116
$ws.Range("G141").Value = 'https://www.olx.pl/d/oferta/kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger.html'
$ws.Range("H141").Value = 'kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger'

# ---- Row 142 ----
$ws.Range("A142").Value = '2026-02-22 12:40:56'
$ws.Range("B142").Value = 'poqui'
$ws.Range("C142").Value = 'Przytulny pokój blisko Politechniki – ul. Przytulna'
$ws.Range("D142").Value = 549
$ws.Range("E142").Value = '10.10.2025'
$ws.Range("F142").Value = 135
$ws.Range("G142").Value = 'https://www.olx.pl/d/oferta/przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz.html'
$ws.Range("H142").Value = 'przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz'

# ---- Row 143 ----
$ws.Range("A143").Value = '2026-02-22 12:40:56'
$ws.Range("B143").Value = 'pokojewlublinie'
$ws.Range("C143").Value = 'WOLNY OD ZARAZ! Pokój jedynka, ul. Romanowskiego 58'
$ws.Range("D143").Value = 0
$ws.Range("E143").Value = '11.08.2025'
$ws.Range("F143").Value = 195
$ws.Range("G143").Value = 'https://www.olx.pl/d/oferta/wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm.html'
$ws.Range("H143").Value = 'wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm'

# ---- Row 144 ----
$ws.Range("A144").Value = '2026-02-22 12:40:56'
$ws.Range("B144").Value = 'pokojewlublinie'
$ws.Range("C144").Value = 'WOLNY OD ZARAZ! Super lokalizacja, blisko centrum, ul. Paganiniego 12'
$ws.Range("D144").Value = 12640
$ws.Range("E144").Value = '19.01.2026'
$ws.Range("F144").Value = 33
$ws.Range("G144").Value = 'https://www.olx.pl/d/oferta/wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc.html'
$ws.Range("H144").Value = 'wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc'

# ---- Row 145 ----
$ws.Range("A145").Value = '2026-02-22 12:40:56'
$ws.Range("B145").Value = 'dawnypatron'
$ws.Range("C145").Value = 'Ładny pokój jednoosobowy. Wynajmę duży pokój w centrum. ul Niecała 4.'
$ws.Range("D145").Value = 730
$ws.Range("E145").Value = '20.09.2024'
$ws.Range("F145").Value = 520
$ws.Range("G145").Value = 'https://www.olx.pl/d/oferta/ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM.html'
$ws.Range("H145").Value = 'ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM'

# ---- Row 146 ----
$ws.Range("A146").Value = '2026-02-22 12:40:56'
$ws.Range("B146").Value = 'dawnypatron'
$ws.Range("C146").Value = 'Mam do wynajęcia pokój dla os. pracującej lub studenta. Narutowicza 14'
$ws.Range("D146").Value = 14690
$ws.Range("E146").Value = '05.12.2025'
$ws.Range("F146").Value = 79
$ws.Range("G146").Value = 'https://www.olx.pl/d/oferta/mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv.html'
$ws.Range("H146").Value = 'mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv'

# ---- Formatting: copy the number/alignment/font formats already used
#      for the data columns from the existing rows, cell by cell, so the
#      new rows look identical to the old ones.
#      Row 7  -> A/C left-aligned text; D/E/F centered, black ("<=60 days")
#      Row 8  -> F centered, red font (">60 days")
$ws.Range("A7").Copy()
$ws.Range("A139").PasteSpecial(-4122)
$ws.Range("C7").Copy()
$ws.Range("C139").PasteSpecial(-4122)
$ws.Range("D7").Copy()
$ws.Range("D139").PasteSpecial(-4122)
$ws.Range("E7").Copy()
$ws.Range("E139").PasteSpecial(-4122)
$ws.Range("F8").Copy()
$ws.Range("F139").PasteSpecial(-4122)
$ws.Range("A7").Copy()
$ws.Range("A140").PasteSpecial(-4122)
$ws.Range("C7").Copy()
$ws.Range("C140").PasteSpecial(-4122)
$ws.Range("D7").Copy()
$ws.Range("D140").PasteSpecial(-4122)
$ws.Range("E7").Copy()
$ws.Range("E140").PasteSpecial(-4122)
$ws.Range("F7").Copy()
$ws.Range("F140").PasteSpecial(-4122)
$ws.Range("A7").Copy()
$ws.Range("A141").PasteSpecial(-4122)
$ws.Range("C7").Copy()
$ws.Range("C141").PasteSpecial(-4122)
$ws.Range("D7").Copy()
$ws.Range("D141").PasteSpecial(-4122)
$ws.Range("E7").Copy()
$ws.Range("E141").PasteSpecial(-4122)
$ws.Range("F8").Copy()
$ws.Range("F141").PasteSpecial(-4122)
$ws.Range("A7").Copy()
$ws.Range("A142").PasteSpecial(-4122)
$ws.Range("C7").Copy()
$ws.Range("C142").PasteSpecial(-4122)
$ws.Range("D7").Copy()
$ws.Range("D142").PasteSpecial(-4122)
$ws.Range("E7").Copy()
$ws.Range("E142").PasteSpecial(-4122)
$ws.Range("F8").Copy()
$ws.Range("F142").PasteSpecial(-4122)
$ws.Range("A7").Copy()
$ws.Range("A143").PasteSpecial(-4122)
$ws.Range("C7").Copy()
$ws.Range("C143").PasteSpecial(-4122)
$ws.Range("D7").Copy()
$ws.Range("D143").PasteSpecial(-4122)
$ws.Range("E7").Copy()
$ws.Range("E143").PasteSpecial(-4122)
$ws.Range("F8").Copy()
$ws.Range("F143").PasteSpecial(-4122)
$ws.Range("A7").Copy()
$ws.Range("A144").PasteSpecial(-4122)
$ws.Range("C7").Copy()
$ws.Range("C144").PasteSpecial(-4122)
$ws.Range("D7").Copy()
$ws.Range("D144").PasteSpecial(-4122)
$ws.Range("E7").Copy()
$ws.Range("E144").PasteSpecial(-4122)
$ws.Range("F7").Copy()
$ws.Range("F144").PasteSpecial(-4122)
$ws.Range("A7").Copy()
$ws.Range("A145").PasteSpecial(-4122)
$ws.Range("C7").Copy()
$ws.Range("C145").PasteSpecial(-4122)
$ws.Range("D7").Copy()
$ws.Range("D145").PasteSpecial(-4122)
$ws.Range("E7").Copy()
$ws.Range("E145").PasteSpecial(-4122)
$ws.Range("F8").Copy()
$ws.Range("F145").PasteSpecial(-4122)
$ws.Range("A7").Copy()
$ws.Range("A146").PasteSpecial(-4122)
$ws.Range("C7").Copy()
$ws.Range("C146").PasteSpecial(-4122)
$ws.Range("D7").Copy()
$ws.Range("D146").PasteSpecial(-4122)
$ws.Range("E7").Copy()
$ws.Range("E146").PasteSpecial(-4122)
$ws.Range("F8").Copy()
$ws.Range("F146").PasteSpecial(-4122)

$excel.CutCopyMode = $false

